$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-12-11 Wednesday" "2024-12-12 Thursday"

Replace-Text "11÷6=" "53÷8="
Replace-Text "64÷7=" "74÷6="
Replace-Text "55÷3=" "96÷9="
Replace-Text "46÷2=" "36÷5="
Replace-Text "59÷7=" "87÷2="

Replace-Text "35÷6=" "90÷7="
Replace-Text "12÷6=" "86÷5="
Replace-Text "55÷7=" "59÷5="
Replace-Text "88÷2=" "43÷5="
Replace-Text "69÷4=" "36÷3="

Replace-Text "56÷4=" "49÷9="
Replace-Text "65÷5=" "92÷6="
Replace-Text "91÷6=" "65÷5="
Replace-Text "46÷9=" "75÷5="
Replace-Text "97÷8=" "89÷7="

Replace-Text "95÷5=" "30÷6="
Replace-Text "51÷8=" "73÷9="
Replace-Text "25÷5=" "95÷3="
Replace-Text "74÷7=" "15÷6="
Replace-Text "40÷6=" "55÷4="

Replace-Text "93÷8=" "81÷5="
Replace-Text "26÷7=" "69÷8="
Replace-Text "45÷8=" "79÷7="
Replace-Text "31÷5=" "23÷5="
Replace-Text "20÷7=" "78÷3="
